# Add a new "additional_accounts" row (covid / other misc accounts) to the
# historical_accounts sheet, as the 13th data row (row index 13, since row 1
# is the header and rows 2-12 hold the existing accounts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account label in column A -> becomes a brand-new shared string.
$ws.Range("A13").Value = "additional_accounts"

# Columns B ("e_trend") and C ("e_cycle") hold the text markers "1.0"/"0.0"
# on every existing data row (shared strings, not numbers). Copy them from
# an existing row instead of typing the values, so they land as the same
# shared-string-typed text (no numeric coercion, no new cell style).
$ws.Range("B2").Copy($ws.Range("B13"))
$ws.Range("C2").Copy($ws.Range("C13"))

# Yearly columns D:S (2006-2021) are all zero for this new account.
$ws.Range("D13:S13").Value = 0

# Match the author's final selection/cursor position recorded in the diff.
$ws.Range("G9").Select() | Out-Null
